# Inclusão de uma nova pesquisa feita e início dos trabalhos dos artigos dentro dessa pesquisa

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cenarios_Pesquisa_WeakSignals")

# ---------------------------------------------------------------------------
# New row 16: a new "pesquisa" search string (col B, rich-text) together with
# its result count (col C).
# ---------------------------------------------------------------------------
$ws.Rows.Item(16).RowHeight = 69.75

$cell = $ws.Cells.Item(16, 2)

# Start from the same base formatting as the other search-string rows
# (B13:B15 - gray 18pt Arial, wrapped) and only tweak the vertical alignment,
# instead of building the font up from scratch (keeps the style table tidy).
$ws.Range("B15").Copy()
$cell.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$cell.Value = '( ( ( TITLE-ABS-KEY ( "foresight" )  AND  ( "competitive intelligence"  OR  "big data"  OR  "forecast" ) )  AND  PUBYEAR  >  2009 ) ) '
$cell.VerticalAlignment = -4108   # xlVAlignCenter

# Colour the literal search terms and the year black, and explicitly restate
# the boilerplate ("AND", "OR", parentheses, ...) in gray - mirrors the
# rich-text runs used by every other search-string cell in this sheet, where
# every run (bar the very first) carries an explicit <rPr>.
$black = 0          # RGB(0,0,0)
$gray  = 9868950    # RGB(150,150,150) = 0x969696

$cell.Characters(23, 11).Font.Color = $black     # "foresight"
$cell.Characters(34, 11).Font.Color = $gray      #  )  AND  (
$cell.Characters(45, 26).Font.Color = $black     # "competitive intelligence"
$cell.Characters(71, 6).Font.Color = $gray       #   OR
$cell.Characters(77, 10).Font.Color = $black     # "big data"
$cell.Characters(87, 6).Font.Color = $gray       #   OR
$cell.Characters(93, 10).Font.Color = $black     # "forecast"
$cell.Characters(103, 23).Font.Color = $gray     # ) )  AND  PUBYEAR  >
$cell.Characters(126, 4).Font.Color = $black     # 2009
$cell.Characters(130, 5).Font.Color = $gray      # ) )

$ws.Cells.Item(16, 3).Value = 708

# ---------------------------------------------------------------------------
# View / navigation state left by the author after the edit.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("C17").Select()
